# Aggiunto Doc Riuso Codice
# Adds a new row (row 28) to the time-tracking table on Foglio1, mirroring
# the pattern of the existing rows (Persona, Progetto, Attività, Data, Tempo).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Copy the date cell's style from the row above so the new date cell keeps
# the same existing numeric date format (instead of creating a new one).
$ws.Range("H27").Copy($ws.Range("H28"))

$ws.Range("E28").Value = "LucaP"
$ws.Range("F28").Value = "GDPR"
$ws.Range("G28").Value = "Form e PHP"
$ws.Range("H28").Value = 43525
$ws.Range("I28").Value = 30

$ws.Range("E29").Select()
